$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "0_5_T1_SP_GRU_EG0_3101_03_F-P-001 - Wand_BA_1_Beton@Erde"
$ws.Range("A3").Value = "0_5_T1_SP_GRU_EG0_3101_03_F-P-001 - Wand_BA_1_Beton@Schalen"
$ws.Range("A4").Value = "0_5_T1_SP_GRU_EG0_3101_03_F-P-001 - Wand_BA_1_Beton@Vorbereitung"
$ws.Range("D3").Value = 45090
